$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull data / push all data / mean calculation - update dSF (column F) values
$ws.Range("F2").Value = 8
$ws.Range("F5").Value = 8
$ws.Range("F7").Value = -5
$ws.Range("F12").Value = -5
$ws.Range("F13").Value = -7
$ws.Range("F18").Value = -2
$ws.Range("F29").Value = -5
$ws.Range("F30").Value = -6
$ws.Range("F36").Value = -2
